$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 51
$ws.Range("I2").Value = 181
$ws.Range("J2").Value = 685
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 170
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 132
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 76
$ws.Range("T2").Value = 140
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 1068
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1086
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 4
